$wb = $excel.ActiveWorkbook

# --- Sheet1: rename and restructure ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Capture existing "ds" (order week) dates from column A before inserting a column
$oldA = @()
for ($r = 2; $r -le 16; $r++) {
    $oldA += $ws1.Cells.Item($r, 1).Value2
}

# Insert a new column C ("Order Week"); old column C (PO_Requested_Qty) shifts to D
$ws1.Columns.Item(3).Insert()

$ws1.Cells.Item(1, 3).Value = "Order Week"

for ($r = 2; $r -le 16; $r++) {
    $ws1.Cells.Item($r, 3).Value = $oldA[$r - 2]
    $ws1.Cells.Item($r, 1).Value = $oldA[$r - 2] + 6
}

$ws1.Range("C2:C16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Add the three new sheets, in order, after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Prediction Info"

# Match outline / page setup markers and page margins used on the original sheet
foreach ($sheet in @($ws2, $ws3, $ws4)) {
    $sheet.Outline.SummaryRow = 1
    $sheet.Outline.SummaryColumn = 1
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# --- Sheet2: Weekly Growth ---
$ws1.Range("A1").Copy($ws2.Range("A1:C1"))
$ws2.Cells.Item(1, 1).Value = "ds"
$ws2.Cells.Item(1, 2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1, 3).Value = "Growth%"

# --- Sheet3: Volume Insights ---
$ws1.Range("A1").Copy($ws3.Range("A1:D1"))
$ws3.Cells.Item(1, 1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1, 2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1, 3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1, 4).Value = "Min_PO_Quantity"

$ws3.Cells.Item(2, 1).Value = 0
$ws3.Cells.Item(2, 2).Value = 0
$ws3.Cells.Item(2, 3).Value = 0
$ws3.Cells.Item(2, 4).Value = 0

# --- Sheet4: Prediction Info ---
$ws1.Range("A1").Copy($ws4.Range("A1"))
$ws4.Cells.Item(1, 1).Value = "Predicted_Next_Week_PO_Quantity"

$ws4.Cells.Item(2, 1).Value = 0

# Make sure the first sheet remains the active one
$ws1.Activate()
